{"js": "// Update the date heading and regenerate the division-practice answer grid.\n//\n// The document is a single title paragraph (\"YYYY-MM-DD Weekday\") followed\n// by one 5-column table whose rows alternate between a row of 5 filled-in\n// answer cells and 3 blank spacer rows. We only ever touch the *text* of\n// existing runs/cells - no rows, columns or paragraphs are added or removed.\n\n// --- 1. Update the date heading -------------------------------------------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst heading = paragraphs.items[0];\nheading.load(\"text\");\nawait context.sync();\n\nif (heading.text.replace(/\\r$/, \"\") === \"2025-06-01 Sunday\") {\n  heading.insertText(\"2025-06-02 Monday\", Word.InsertLocation.replace);\n}\n\n// --- 2. Update the answer table --------------------------------------------------\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, row by row, in the same 20x5 shape as the existing table\n// (the 3 blank spacer rows after each answer row stay blank; only the\n// five \"answer\" rows actually carry text).\ntable.values = [\n  [\"18\u00f74=4, 2\", \"62\u00f73=20, 2\", \"85\u00f78=10, 5\", \"80\u00f76=13, 2\", \"37\u00f73=12, 1\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"98\u00f76=16, 2\", \"85\u00f73=28, 1\", \"75\u00f74=18, 3\", \"82\u00f73=27, 1\", \"69\u00f73=23, 0\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"15\u00f72=7, 1\", \"52\u00f74=13, 0\", \"21\u00f74=5, 1\", \"74\u00f72=37, 0\", \"81\u00f78=10, 1\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"39\u00f77=5, 4\", \"13\u00f76=2, 1\", \"23\u00f79=2, 5\", \"56\u00f74=14, 0\", \"94\u00f79=10, 4\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"17\u00f76=2, 5\", \"69\u00f76=11, 3\", \"93\u00f75=18, 3\", \"35\u00f75=7, 0\", \"28\u00f72=14, 0\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n];\n\nawait context.sync();\n", "ps1": "# Update the date heading and regenerate the division-practice answer grid.\n#\n# The document is a single title paragraph (\"YYYY-MM-DD Weekday\") followed\n# by one 5-column table whose rows alternate between a row of 5 filled-in\n# answer cells and 3 blank spacer rows. We only ever touch the *text* of\n# existing runs/cells - no rows, columns or paragraphs are added or removed.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date heading -------------------------------------------------\n$heading = $d.Paragraphs.Item(1)\n$headingText = $heading.Range.Text -replace \"[\\r\\x07]+$\", \"\"\nif ($headingText -eq \"2025-06-01 Sunday\") {\n    $rng = $heading.Range\n    $rng.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; keep the paragraph mark\n    $rng.Text = \"2025-06-02 Monday\"\n}\n\n# --- 2. Update the answer table --------------------------------------------------\n$t = $d.Tables.Item(1)\n\n# New values, row by row, in the same 20x5 shape as the existing table\n# (the 3 blank spacer rows after each answer row stay blank; only the\n# five \"answer\" rows actually carry text).\n$newValues = @(\n    @(\"18\u00f74=4, 2\", \"62\u00f73=20, 2\", \"85\u00f78=10, 5\", \"80\u00f76=13, 2\", \"37\u00f73=12, 1\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"98\u00f76=16, 2\", \"85\u00f73=28, 1\", \"75\u00f74=18, 3\", \"82\u00f73=27, 1\", \"69\u00f73=23, 0\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"15\u00f72=7, 1\", \"52\u00f74=13, 0\", \"21\u00f74=5, 1\", \"74\u00f72=37, 0\", \"81\u00f78=10, 1\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"39\u00f77=5, 4\", \"13\u00f76=2, 1\", \"23\u00f79=2, 5\", \"56\u00f74=14, 0\", \"94\u00f79=10, 4\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"17\u00f76=2, 5\", \"69\u00f76=11, 3\", \"93\u00f75=18, 3\", \"35\u00f75=7, 0\", \"28\u00f72=14, 0\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\"),\n    @(\"\", \"\", \"\", \"\", \"\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $newText = $rowValues[$c - 1]\n        if ($newText -ne \"\") {\n            $cell = $t.Cell($r, $c)\n            $cellRng = $cell.Range\n            $cellRng.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; keep the cell mark\n            $cellRng.Text = $newText\n        }\n    }\n}\n"}
